$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets carry the same table and need the same
# F-column ("想去人数") updates for rows 4, 7 and 10.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 1562
    $ws.Range("F7").Value = 387
    $ws.Range("F10").Value = 415
}
